$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "73-32=", "78-28=", "30+51=", "58+19=", "65-12=", "76-26=", "27+13=", "69-49=", "77-43=", "35+23=",
    "12+45=", "45+26=", "24-2=", "23+6=", "63-61=", "35-12=", "74-41=", "51-49=", "91-77=", "22+21=",
    "12+63=", "32+8=", "11+73=", "25+21=", "41+16=", "65-44=", "75-25=", "51+33=", "75-69=", "63+34=",
    "32-16=", "7+56=", "2+43=", "75-64=", "97-8=", "48+12=", "90-33=", "93-66=", "42-5=", "81+11=",
    "41-33=", "55+15=", "52+28=", "61-33=", "65+19=", "77-52=", "3+83=", "9+27=", "89-23=", "72-16=",
    "68-1=", "97-71=", "39+27=", "1+29=", "38+56=", "39+59=", "44-33=", "46+47=", "98-63=", "56-6=",
    "99-77=", "3+71=", "65+2=", "88-13=", "2+44=", "56-40=", "4+16=", "27+22=", "46+12=", "78+8=",
    "90-76=", "77+11=", "41-23=", "85-61=", "44-38=", "2+24=", "2+14=", "83-1=", "98-84=", "84+15=",
    "73+24=", "2+11=", "34-17=", "52-5=", "51-30=", "76+14=", "23+5=", "38+1=", "48+11=", "6+6=",
    "20+22=", "92-27=", "88-16=", "85-12=", "70-35=", "81-61=", "81-74=", "2+40=", "16+53=", "38-18="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
